# Update the division-fact answers in the practice table.
# Each table cell holds one line of text like "55÷3=18, 1" which is
# replaced with a new fact/answer pair, per the commit's regenerated output.

$d = $word.ActiveDocument

# --- Straightforward one-to-one text replacements (each old string is
#     unique in the document, so a simple Find/Replace is unambiguous). ---
$replacements = @(
    @("55÷3=18, 1", "11÷7=1, 4"),
    @("66÷2=33, 0", "16÷9=1, 7"),
    @("75÷8=9, 3", "16÷7=2, 2"),
    @("14÷8=1, 6", "56÷7=8, 0"),
    @("35÷6=5, 5", "95÷7=13, 4"),
    @("55÷8=6, 7", "22÷2=11, 0"),
    @("77÷9=8, 5", "10÷6=1, 4"),
    @("61÷7=8, 5", "63÷6=10, 3"),
    @("99÷8=12, 3", "35÷5=7, 0"),
    @("32÷5=6, 2", "84÷8=10, 4"),
    @("77÷5=15, 2", "57÷2=28, 1"),
    @("97÷8=12, 1", "85÷7=12, 1"),
    @("47÷7=6, 5", "95÷4=23, 3"),
    @("93÷8=11, 5", "82÷7=11, 5"),
    @("64÷8=8, 0", "73÷9=8, 1"),
    @("85÷2=42, 1", "71÷4=17, 3"),
    @("75÷2=37, 1", "82÷8=10, 2"),
    @("77÷3=25, 2", "71÷9=7, 8"),
    @("43÷8=5, 3", "73÷7=10, 3"),
    @("18÷7=2, 4", "90÷7=12, 6")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Last row of the table: every cell's value changes, and the new value
#     of one cell equals the *old* value of a neighboring cell
#     ("54÷9=6, 0" -> "71÷6=11, 5", while the original "71÷6=11, 5" cell
#     becomes "88÷5=17, 3"). A plain document-wide Find/Replace would be
#     ambiguous here, so address these cells positionally instead. ---
$table = $d.Tables.Item(1)
$lastRow = 17

$table.Cell($lastRow, 1).Range.Text = "35÷2=17, 1"
$table.Cell($lastRow, 2).Range.Text = "48÷3=16, 0"
$table.Cell($lastRow, 3).Range.Text = "71÷6=11, 5"
$table.Cell($lastRow, 4).Range.Text = "88÷5=17, 3"
$table.Cell($lastRow, 5).Range.Text = "50÷6=8, 2"
